$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("People")

# Fill down column B ("quote" helper column) from row 35 into rows 36:46,
# and rebuild the E (Name link) / I (Institution link) helper-formulas for
# those same rows, mirroring the pattern already used in rows 2:35.
for ($r = 36; $r -le 46; $r++) {
    $ws.Range("B$r").Value = """"

    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Formula = '="<a href="&B' + $r + '&C' + $r + '&B' + $r + '&">"&D' + $r + '&"</a>"'

    $ws.Range("I$r").Formula = '="<a href="&B' + $r + '&G' + $r + '&B' + $r + '&">"&H' + $r + '&"</a>"'
}

# The sheet had scrolled so column E was the left-most visible column and
# H52 was selected; reset the scroll position and move the selection to F53.
$w = $wb.Windows.Item(1)
$w.ScrollColumn = 1
$w.ScrollRow = 1
$ws.Range("F53").Select()
